$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K11").Value = 0.1683237681281231
$ws.Range("J12").Value = 0.1722916656412322
$ws.Range("I13").Value = 0.3079317558114735
$ws.Range("H14").Value = 0.06712557395580883
$ws.Range("G15").Value = 0.02179435870371246
$ws.Range("F16").Value = -0.04506706323234141
$ws.Range("E17").Value = -0.07465326558905801
$ws.Range("D18").Value = -0.0928039223186989
$ws.Range("C19").Value = -0.1108357465673982
$ws.Range("B20").Value = -0.1624199859130616
